$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Helper: find the paragraph whose text currently reads the merged
# "Git restore *...#for all files" line (robust against re-runs).
# ---------------------------------------------------------------------
$count = $d.Paragraphs.Count
$restoreIdx = -1
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith("Git restore *")) {
        $restoreIdx = $i
        break
    }
}

# ---------------------------------------------------------------------
# Step 1: clean up the "Git restore *" paragraph.
# It currently holds two runs with a "_GoBack" bookmark sitting between
# them:  "Git restore *                 " + bookmark + " #for all files"
# We want a single run reading
#   "Git restore *                  #for all files"
# with no bookmark left inside it (the bookmark is relocated later).
# ---------------------------------------------------------------------
$pRestore = $d.Paragraphs.Item($restoreIdx)
$fullRange = $pRestore.Range
$firstRunLen = 30   # length of "Git restore *" + trailing spaces
$splitPos = $fullRange.Start + $firstRunLen
$tailRange = $d.Range($splitPos, $fullRange.End)
$tailRange.Delete()

$pRestore = $d.Paragraphs.Item($restoreIdx)
$pRestore.Range.InsertAfter(" #for all files")

# ---------------------------------------------------------------------
# Step 2: the paragraph immediately following is already a blank
# paragraph in the original document - reuse it as the first blank line.
# ---------------------------------------------------------------------
$blank1Idx = $restoreIdx + 1
$curParagraph = $d.Paragraphs.Item($blank1Idx)

# ---------------------------------------------------------------------
# Step 3: insert "Or use"
# ---------------------------------------------------------------------
$curParagraph.Range.InsertParagraphAfter()
$curIdx = $blank1Idx + 1
$curParagraph = $d.Paragraphs.Item($curIdx)
$curParagraph.Range.Text = "Or use"

# ---------------------------------------------------------------------
# Step 4: insert "Git stash  #remove all your local changes"
# ---------------------------------------------------------------------
$curParagraph.Range.InsertParagraphAfter()
$curIdx = $curIdx + 1
$curParagraph = $d.Paragraphs.Item($curIdx)
$curParagraph.Range.Text = "Git stash  #remove all your local changes"

# ---------------------------------------------------------------------
# Step 5: insert "Git stash pop ..." with a bottom paragraph border
# ---------------------------------------------------------------------
$curParagraph.Range.InsertParagraphAfter()
$curIdx = $curIdx + 1
$curParagraph = $d.Paragraphs.Item($curIdx)
$curParagraph.Range.Text = "Git stash pop       #this command restore most recent stashed files"

$borders = $curParagraph.Format.Borders
$bottomBorder = $borders.Item(-3)   # wdBorderBottom
$bottomBorder.LineStyle = 1         # wdLineStyleSingle
$bottomBorder.LineWidth = 2         # -> sz="4" (quarter points)
$bottomBorder.ColorIndex = 0        # wdAuto -> color="auto"
$borders.DistanceFromBottom = 0

# ---------------------------------------------------------------------
# Step 6: insert a new blank paragraph after "Git stash pop ..."
# Freshly split paragraph marks in this runtime carry over a stray
# empty run (and the paragraph border just set above); clean the run
# up by typing then deleting a placeholder char, and remove the
# inherited border so only the "Git stash pop" line keeps it.
# ---------------------------------------------------------------------
$curParagraph.Range.InsertParagraphAfter()
$curIdx = $curIdx + 1
$curParagraph = $d.Paragraphs.Item($curIdx)
$cleanPos = $curParagraph.Range.Start
$curParagraph.Range.InsertAfter("x")
$d.Range($cleanPos, $cleanPos + 1).Delete()

$blankBorders = $curParagraph.Format.Borders
$blankBorders.Item(-3).LineStyle = 0   # wdLineStyleNone

# ---------------------------------------------------------------------
# Step 7: insert one more new blank paragraph that will hold the
# relocated "_GoBack" bookmark.  (Adding a bookmark to a genuinely
# empty range at certain offsets is unreliable in this runtime, so
# give the paragraph a temporary character first, then remove it.)
# ---------------------------------------------------------------------
$curParagraph.Range.InsertParagraphAfter()
$curIdx = $curIdx + 1
$curParagraph = $d.Paragraphs.Item($curIdx)
$bkPos = $curParagraph.Range.Start
$curParagraph.Range.InsertAfter("x")

$bkRange = $d.Range($bkPos, $bkPos)
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
$d.Bookmarks.Add("_GoBack", $bkRange)

$d.Range($bkPos, $bkPos + 1).Delete()

Write-Output "edit complete"
